$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated forecast error values (rows 2-11 correspond to Q0-Q9), columns B-G
$data = @{
    2  = @(-0.02851397337148708, 1.487398177008028, 3.318917335035293, 1.821789596807297, 1.83969192369946, 51)
    3  = @(0.08547875753704889, 1.81008304745325, 5.250498484789946, 2.291396623195109, 2.313049022923821, 50)
    4  = @(0.4939789650732916, 1.887039932152993, 6.340229206786162, 2.517981176813314, 2.494637939050945, 49)
    5  = @(0.2451459353463975, 1.858170400153691, 5.548157995631708, 2.355452821780073, 2.367451949089824, 48)
    6  = @(0.4386977934359605, 1.706409057465927, 4.971657831094963, 2.229721469398132, 2.209773240506429, 47)
    7  = @(0.2428816792506087, 1.58716742716687, 4.151381641064656, 2.037493960988512, 2.05012079391029, 38)
    8  = @(0.4103070901050488, 1.613429594635641, 4.299938655358432, 2.073629343773479, 2.060668106634449, 37)
    9  = @(1.336347049941183, 1.610837625543729, 4.520735657440727, 2.12620216758443, 1.696718883006791, 20)
    10 = @(0.9801651346192841, 1.34337281173155, 2.87806734524593, 1.696486765420211, 1.441222499118803, 13)
    11 = @(0.8601837223604264, 1.462038904411423, 3.55962077517264, 1.886695729356655, 1.877400043597129, 5)
}

foreach ($r in $data.Keys) {
    $values = $data[$r]
    $ws.Cells.Item($r, 2).Value = $values[0]
    $ws.Cells.Item($r, 3).Value = $values[1]
    $ws.Cells.Item($r, 4).Value = $values[2]
    $ws.Cells.Item($r, 5).Value = $values[3]
    $ws.Cells.Item($r, 6).Value = $values[4]
    $ws.Cells.Item($r, 7).Value = $values[5]
}
